$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.954.26"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "2.362.23"
$ws.Range("E3").Value = "  +1.91%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'0.662"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.55%  "
$ws.Range("D6").Value = "'235.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.28%  "
$ws.Range("D7").Value = "'72.58"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +11.45%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +21.40%  "
$ws.Range("D10").Value = "'0.0986"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.52%  "
$ws.Range("D11").Value = "'28.25"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.85%  "
$ws.Range("D12").Value = "2.713.27"
$ws.Range("E12").Value = "  +2.82%  "
$ws.Range("E13").Value = "  +1.62%  "
$ws.Range("D14").Value = "'16.84"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.59%  "
$ws.Range("D15").Value = "'6.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +8.18%  "
$ws.Range("D16").Value = "'0.882"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.68%  "
$ws.Range("D17").Value = "2.365.62"
$ws.Range("E17").Value = "  +2.59%  "
$ws.Range("D18").Value = "43.899.34"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("E19").Value = "  +3.02%  "
$ws.Range("D20").Value = "'76.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.67%  "
$ws.Range("D21").Value = "'6.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.49%  "
$ws.Range("D22").Value = "'251.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.78%  "
$ws.Range("E23").Value = "  +1.80%  "
$ws.Range("D24").Value = "'0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("E25").Value = "  +2.12%  "
$ws.Range("D26").Value = "'10.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.70%  "
$ws.Range("D27").Value = "'2.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("D28").Value = "'22.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.86%  "
$ws.Range("D29").Value = "'172.89"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.88%  "
$ws.Range("E30").Value = "  +8.61%  "
$ws.Range("D31").Value = "'0.132"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.34%  "
$ws.Range("E32").Value = "  +4.46%  "
$ws.Range("D33").Value = "'5.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.23%  "
$ws.Range("E34").Value = "  +3.64%  "
$ws.Range("D35").Value = "'5.15"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.49%  "
$ws.Range("D36").Value = "'3.75"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.90%  "
$ws.Range("D37").Value = "'2.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.31%  "
$ws.Range("D38").Value = "'6.42"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.72%  "
$ws.Range("E39").Value = "  +5.47%  "
$ws.Range("D40").Value = "'19.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.48%  "
$ws.Range("B41").Value = "BinanceUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'8.94"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").Value = "'1.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.09%  "
$ws.Range("E44").Value = "  +1.95%  "
$ws.Range("D45").Value = "'1.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.25%  "
$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").Value = "'4.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.43%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'98.09"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("E48").Value = "  +12.67%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "1.438.88"
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'2.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.59%  "
$ws.Range("D51").Value = "2.586.54"
$ws.Range("E51").Value = "  +2.60%  "
